$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.674.82'
$ws.Range("E2").Value = '  -0.25%  '

$ws.Range("D3").Value = '2.679.31'
$ws.Range("E3").Value = '  -0.56%  '

$ws.Range("E4").Value = '  -0.07%  '

$ws.Range("D5").Value = '''601.24'
$ws.Range("E5").Value = '  -1.27%  '

$ws.Range("D6").Value = '''157.09'
$ws.Range("E6").Value = '  -0.22%  '

$ws.Range("D8").Value = '''0.624'
$ws.Range("E8").Value = '  +6.04%  '

$ws.Range("D9").Value = '''0.130'
$ws.Range("E9").Value = '  +4.54%  '

$ws.Range("D10").Value = '''0.402'
$ws.Range("E10").Value = '  -0.42%  '

$ws.Range("D11").Value = '''5.86'
$ws.Range("E11").Value = '  -2.86%  '

$ws.Range("E12").Value = '  -0.25%  '

$ws.Range("D13").Value = '''29.39'
$ws.Range("E13").Value = '  -2.77%  '

$ws.Range("D14").Value = '''0.0000199'
$ws.Range("E14").Value = '  -1.17%  '

$ws.Range("D15").Value = '3.156.91'
$ws.Range("E15").Value = '  -0.71%  '

$ws.Range("D16").Value = '65.522.72'
$ws.Range("E16").Value = '  -0.25%  '

$ws.Range("D17").Value = '2.675.39'
$ws.Range("E17").Value = '  -0.84%  '

$ws.Range("D18").Value = '''12.81'
$ws.Range("E18").Value = '  +1.10%  '

$ws.Range("E19").Value = '  -1.75%  '

$ws.Range("D20").Value = '''7.59'
$ws.Range("E20").Value = '  -0.28%  '

$ws.Range("D21").Value = '''352.01'
$ws.Range("E21").Value = '  -2.09%  '

$ws.Range("E22").Value = '  -0.02%  '

$ws.Range("D23").Value = '''69.60'
$ws.Range("E23").Value = '  -0.83%  '

$ws.Range("D24").Value = '''0.0000112'
$ws.Range("E24").Value = '  +6.05%  '

$ws.Range("D25").Value = '''9.69'
$ws.Range("E25").Value = '  -0.60%  '

$ws.Range("E26").Value = '  +0.75%  '

$ws.Range("B27").Value = 'Kaspa'
$ws.Range("C27").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D27").Value = '''0.167'
$ws.Range("E27").Value = '  -0.76%  '

$ws.Range("B28").Value = 'Fetch.AI'
$ws.Range("C28").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D28").Value = '''1.60'
$ws.Range("E28").Value = '  -5.46%  '

$ws.Range("D29").Value = '''8.12'
$ws.Range("E29").Value = '  -1.01%  '

$ws.Range("E30").Value = '  -0.22%  '

$ws.Range("D31").Value = '''534.93'
$ws.Range("E31").Value = '  +0.46%  '

$ws.Range("D32").Value = '''2.15'
$ws.Range("E32").Value = '  -2.32%  '

$ws.Range("D33").Value = '''1.75'
$ws.Range("E33").Value = '  -2.31%  '

$ws.Range("D34").Value = '''6.47'
$ws.Range("E34").Value = '  -2.38%  '

$ws.Range("D35").Value = '''5.51'
$ws.Range("E35").Value = '  +1.15%  '

$ws.Range("D36").Value = '''0.425'
$ws.Range("E36").Value = '  -1.78%  '

$ws.Range("D37").Value = '''20.51'
$ws.Range("E37").Value = '  -1.19%  '

$ws.Range("D38").Value = '''1.00'
$ws.Range("E38").Value = '  +0.12%  '

$ws.Range("D39").Value = '''158.18'
$ws.Range("E39").Value = '  -3.15%  '

$ws.Range("E40").Value = '  -2.33%  '

$ws.Range("D42").Value = '''164.69'
$ws.Range("E42").Value = '  -2.72%  '

$ws.Range("D43").Value = '''4.15'
$ws.Range("E43").Value = '  -0.54%  '

$ws.Range("E44").Value = '  +2.66%  '

$ws.Range("D45").Value = '''0.0610'
$ws.Range("E45").Value = '  -0.28%  '

$ws.Range("D46").Value = '''22.86'
$ws.Range("E46").Value = '  -2.76%  '

$ws.Range("D47").Value = '''0.0260'
$ws.Range("E47").Value = '  -2.13%  '

$ws.Range("D48").Value = '''0.642'
$ws.Range("E48").Value = '  -2.37%  '

$ws.Range("D49").Value = '0.0₆0258'
$ws.Range("E49").Value = '  +14.04%  '

$ws.Range("E50").Value = '  +2.47%  '

$ws.Range("D51").Value = '''20.19'
$ws.Range("E51").Value = '  -4.60%  '
